# Major overhaul for weights -- final
# Updates the dag/dhs covariate map on Sheet1 to reflect the new weighting
# scheme: several covariate codes are swapped out for their updated
# equivalents, and the active selection is moved to B2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write new covariate codes/labels in the same order the new shared
# strings were introduced so the resulting sharedStrings table lines up
# with the target workbook.

# C:Hospital Distance -> hlthdist_cont_log_scale_clst (was hlthst_nrst_duration_fctb_clst)
$ws.Range("B12").Value = "hlthdist_cont_log_scale_clst"

# I:Hemoglobin/lowhb_fctb -> I:Anemia/hab57_fctb
$ws.Range("A7").Value = "I:Anemia"
$ws.Range("B7").Value = "hab57_fctb"

# I:Wealth -> wlthrcde_combscor_cont (was wlthrcde_fctb)
$ws.Range("B3").Value = "wlthrcde_combscor_cont"

# I:Education -> hv108_cont_scale (was hv106_fctb)
$ws.Range("B17").Value = "hv108_cont_scale"

# I:Owns Livestock/hv246_fctb -> I:Farmer/farmer_fctb
$ws.Range("A10").Value = "I:Farmer"
$ws.Range("B10").Value = "farmer_fctb"

# Move the active selection to B2
$ws.Range("B2").Select()
